# The deck has two identical textboxes ("3.紀錄遊戲數據") that each hold the
# text as three runs: "3." / "紀錄遊戲" / "數據".
#
# The edit:
#   1. fixes the typo 紀 -> 記 (so the text reads "3.記錄遊戲數據"), and
#   2. splits the runs so that "3." becomes two runs ("3" then "."), and
#      "紀錄遊戲" becomes three runs ("記" then "錄" then "遊戲").
#
# Rather than hard-coding slide/shape indices, find every shape whose text
# starts with "3.紀錄遊戲" and apply the same run split to each of them.

$p = $ppt.ActivePresentation
$needle = "3." + "紀錄遊戲"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }

        $tr = $sh.TextFrame.TextRange
        $text = $tr.Text
        if ($text.IndexOf($needle) -ne 0) { continue }

        # Character positions (1-based) within "3.紀錄遊戲...":
        #   1 = 3   2 = .   3 = 紀(->記)   4 = 錄   5-6 = 遊戲
        $tr.Characters(1, 1).Text = "3"
        $tr.Characters(2, 1).Text = "."
        $tr.Characters(3, 1).Text = "記"
        $tr.Characters(4, 1).Text = "錄"
        $tr.Characters(5, 2).Text = "遊戲"
    }
}
